$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.02
$ws.Range("G2").Value = 2.42
$ws.Range("H2").Value = 2.8
$ws.Range("I2").Value = 3.55
$ws.Range("J2").Value = 3.55
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 1.23
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 3.45
$ws.Range("O2").Value = 1.13
$ws.Range("P2").Value = 2.58
$ws.Range("Q2").Value = 1.4
$ws.Range("R2").Value = 1.71
$ws.Range("S2").Value = 2
$ws.Range("T2").Value = 1.42
$ws.Range("U2").Value = 2.68
$ws.Range("V2").Value = 1.4
$ws.Range("W2").Value = 1.7
$ws.Range("AN2").Value = 55
$ws.Range("N3").Value = 1.02
$ws.Range("P3").Value = 1.28
$ws.Range("N4").Value = 1.02
$ws.Range("P4").Value = 1.28
$ws.Range("R4").Value = 1.21
$ws.Range("F5").Value = 1.8
$ws.Range("G5").Value = 1.81
$ws.Range("H5").Value = 4.3
$ws.Range("K5").Value = 4.7
$ws.Range("W5").Value = 2.22
$ws.Range("AD5").Value = 18
$ws.Range("AF5").Value = 14.5
$ws.Range("AI5").Value = 48
$ws.Range("AK5").Value = 17.5
$ws.Range("AO5").Value = 1000
$ws.Range("J6").Value = 5
$ws.Range("K6").Value = 5.2
$ws.Range("N6").Value = 5.5
$ws.Range("P6").Value = 2.5
$ws.Range("S6").Value = 2.58
$ws.Range("F7").Value = 2.18
$ws.Range("G7").Value = 2.2
$ws.Range("H7").Value = 3.45
$ws.Range("I7").Value = 3.5
$ws.Range("J7").Value = 3.85
$ws.Range("K7").Value = 3.9
$ws.Range("R7").Value = 1.57
$ws.Range("S7").Value = 2.68
$ws.Range("U7").Value = 2.58
$ws.Range("V7").Value = 1.4
$ws.Range("W7").Value = 1.83
$ws.Range("AH7").Value = 15
$ws.Range("AL7").Value = 28
$ws.Range("AO7").Value = 26
$ws.Range("F8").Value = 1.5
$ws.Range("H8").Value = 5.4
$ws.Range("K8").Value = 5.8
$ws.Range("M8").Value = 1.02
$ws.Range("N8").Value = 6.8
$ws.Range("O8").Value = 1.13
$ws.Range("P8").Value = 3.05
$ws.Range("Q8").Value = 1.39
$ws.Range("S8").Value = 1.98
$ws.Range("T8").Value = 1.52
$ws.Range("U8").Value = 2.46
$ws.Range("X8").Value = 40
$ws.Range("Z8").Value = 65
$ws.Range("AB8").Value = 16.5
$ws.Range("AC8").Value = 14.5
$ws.Range("AD8").Value = 25
$ws.Range("AE8").Value = 330
$ws.Range("AG8").Value = 11
$ws.Range("AJ8").Value = 17
$ws.Range("AM8").Value = 65
$ws.Range("AN8").Value = 5.4
$ws.Range("AO8").Value = 50
$ws.Range("J9").Value = 3.65
$ws.Range("F10").Value = 1.52
$ws.Range("G10").Value = 1.58
$ws.Range("K10").Value = 5
$ws.Range("Q10").Value = 1.62
$ws.Range("F11").Value = 3.15
$ws.Range("G11").Value = 3.2
$ws.Range("H11").Value = 2.4
$ws.Range("I11").Value = 2.42
$ws.Range("O11").Value = 1.23
$ws.Range("P11").Value = 2.34
$ws.Range("Q11").Value = 1.72
$ws.Range("T11").Value = 1.61
$ws.Range("V11").Value = 1.7
$ws.Range("W11").Value = 1.45
$ws.Range("Z11").Value = 17
$ws.Range("AI11").Value = 29
$ws.Range("AK11").Value = 30
$ws.Range("AN11").Value = 22
$ws.Range("AO11").Value = 14
$ws.Range("H12").Value = 1.43
$ws.Range("I12").Value = 1.44
$ws.Range("J12").Value = 5.3
$ws.Range("K12").Value = 5.4
$ws.Range("N12").Value = 5.3
$ws.Range("P12").Value = 2.5
$ws.Range("Q12").Value = 1.65
$ws.Range("T12").Value = 1.9
$ws.Range("U12").Value = 2.04
$ws.Range("V12").Value = 3.25
$ws.Range("X12").Value = 24
$ws.Range("Y12").Value = 10
$ws.Range("AC12").Value = 12.5
$ws.Range("AD12").Value = 9.6
$ws.Range("AE12").Value = 14
$ws.Range("AI12").Value = 29
$ws.Range("AJ12").Value = 280
$ws.Range("AN12").Value = 110
$ws.Range("F13").Value = 2.96
$ws.Range("P13").Value = 2.82
$ws.Range("Q13").Value = 1.54
$ws.Range("T13").Value = 1.5
$ws.Range("Z13").Value = 19.5
$ws.Range("AE13").Value = 20
$ws.Range("AL13").Value = 27
$ws.Range("AN13").Value = 15
$ws.Range("H14").Value = 20
$ws.Range("J14").Value = 8.6
$ws.Range("R14").Value = 1.76
$ws.Range("T14").Value = 2.34
$ws.Range("U14").Value = 1.69
$ws.Range("W14").Value = 6
$ws.Range("Y14").Value = 65
$ws.Range("Z14").Value = 1000
$ws.Range("AE14").Value = 460
$ws.Range("AK14").Value = 14.5
$ws.Range("AL14").Value = 50
$ws.Range("O15").Value = 1.27
$ws.Range("P15").Value = 2.16
$ws.Range("Q15").Value = 1.81
$ws.Range("S15").Value = 3.05
$ws.Range("T15").Value = 1.83
$ws.Range("V15").Value = 2.32
$ws.Range("X15").Value = 17
$ws.Range("Y15").Value = 9.800000000000001
$ws.Range("AA15").Value = 16
$ws.Range("AB15").Value = 20
$ws.Range("F16").Value = 1.28
$ws.Range("G16").Value = 1.29
$ws.Range("H16").Value = 12
$ws.Range("I16").Value = 12.5
$ws.Range("J16").Value = 6.8
$ws.Range("K16").Value = 7
$ws.Range("S16").Value = 1.89
$ws.Range("U16").Value = 2.26
$ws.Range("V16").Value = 1.08
$ws.Range("W16").Value = 4.4
$ws.Range("Z16").Value = 130
$ws.Range("AD16").Value = 42
$ws.Range("AI16").Value = 100
$ws.Range("AL16").Value = 26
$ws.Range("M17").Value = 1.09
$ws.Range("N17").Value = 3.45
$ws.Range("O17").Value = 1.38
$ws.Range("P17").Value = 1.82
$ws.Range("Q17").Value = 2.18
$ws.Range("S17").Value = 4
$ws.Range("T17").Value = 1.88
$ws.Range("U17").Value = 2.08
$ws.Range("AN17").Value = 24
$ws.Range("F18").Value = 2.34
$ws.Range("H18").Value = 2.58
$ws.Range("I18").Value = 2.98
$ws.Range("K18").Value = 4.3
$ws.Range("M18").Value = 1.04
$ws.Range("N18").Value = 5.4
$ws.Range("O18").Value = 1.17
$ws.Range("P18").Value = 2.8
$ws.Range("Q18").Value = 1.5
$ws.Range("R18").Value = 1.71
$ws.Range("S18").Value = 2.3
$ws.Range("T18").Value = 1.5
$ws.Range("U18").Value = 2.56
$ws.Range("V18").Value = 1.5
$ws.Range("X18").Value = 1000
$ws.Range("Y18").Value = 18.5
$ws.Range("AA18").Value = 240
$ws.Range("AB18").Value = 20
$ws.Range("AC18").Value = 11
$ws.Range("AG18").Value = 12.5
$ws.Range("AH18").Value = 15.5
$ws.Range("AI18").Value = 34
$ws.Range("AK18").Value = 25
$ws.Range("AM18").Value = 260
$ws.Range("AN18").Value = 14.5
$ws.Range("AO18").Value = 16.5
$ws.Range("G19").Value = 2.8
$ws.Range("I19").Value = 3.25
$ws.Range("P19").Value = 1.83
$ws.Range("Q19").Value = 1.94
$ws.Range("R19").Value = 1.32
$ws.Range("S19").Value = 3.55
$ws.Range("T19").Value = 1.76
$ws.Range("U19").Value = 2.1
$ws.Range("W19").Value = 1.57
$ws.Range("AC19").Value = 7.8
$ws.Range("AE19").Value = 36
$ws.Range("AI19").Value = 130
$ws.Range("AK19").Value = 32
$ws.Range("AM19").Value = 580
$ws.Range("AN19").Value = 27
$ws.Range("AO19").Value = 34
$ws.Range("M20").Value = 1.04
$ws.Range("N20").Value = 1.26
$ws.Range("O20").Value = 1.04
$ws.Range("P20").Value = 1.26
$ws.Range("R20").Value = 1.19
